$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.3558718861209965
$ws1.Range("C2").Value = 0.06510416666666667
$ws1.Range("D2").Value = 0.8928571428571429
$ws1.Range("E2").Value = 0.1213592233009709
$ws1.Range("F2").Value = 0.2520161290322581
$ws1.Range("G2").Value = 0.5996309963099631
$ws1.Range("H2").Value = 0.8013643659711074
$ws1.Range("I2").Value = 25
$ws1.Range("J2").Value = 359
$ws1.Range("K2").Value = 175
$ws1.Range("L2").Value = 3

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.9831460674157303
$ws2.Range("C2").Value = 0.3277153558052435
$ws2.Range("D2").Value = 0.4915730337078651

$ws2.Range("B3").Value = 0.06510416666666667
$ws2.Range("C3").Value = 0.8928571428571429
$ws2.Range("D3").Value = 0.1213592233009709

$ws2.Range("B4").Value = 0.3558718861209965
$ws2.Range("C4").Value = 0.3558718861209965
$ws2.Range("D4").Value = 0.3558718861209965
$ws2.Range("E4").Value = 0.3558718861209965

$ws2.Range("B5").Value = 0.5241251170411985
$ws2.Range("C5").Value = 0.6102862493311931
$ws2.Range("D5").Value = 0.306466128504418

$ws2.Range("B6").Value = 0.937407325029656
$ws2.Range("C6").Value = 0.3558718861209965
$ws2.Range("D6").Value = 0.4731282175310092

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 175
$ws3.Range("C2").Value = 359
$ws3.Range("B3").Value = 3
$ws3.Range("C3").Value = 25
